# ValidateTotals.xlsx edit script
# Implements: add "Input" cell-style headers, shift the right-hand mini-table
# up by one row, and append a second validation block (rows 15-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------------
# Step 0: remember original values of the right-hand table (G3:K6) before we
# start mutating the sheet, since several of them need to move.
# ---------------------------------------------------------------------------
$origG4 = $ws.Range("G4").Value2
$origH4 = $ws.Range("H4").Value2
$origI4 = $ws.Range("I4").Value2
$origG5 = $ws.Range("G5").Value2
$origH5 = $ws.Range("H5").Value2
$origI5 = $ws.Range("I5").Value2

# ---------------------------------------------------------------------------
# Step 1: apply the new (non-bold) "Input" style to the combined header row 2
# (left table A2:D2, already styled "s=3"; right table header moves here).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "DATE"
$ws.Range("A2").Style = "Input"
$ws.Range("B2").Value = "FILE"
$ws.Range("B2").Style = "Input"
$ws.Range("C2").Value = "Count"
$ws.Range("C2").Style = "Input"
$ws.Range("D2").Value = "Difference"
$ws.Range("D2").Style = "Input"

$ws.Range("G2").Value = "DATE"
$ws.Range("G2").Style = "Input"
$ws.Range("H2").Value = "FILE"
$ws.Range("H2").Style = "Input"
$ws.Range("I2").Value = "NVD billed"
$ws.Range("I2").Style = "Input"

# ---------------------------------------------------------------------------
# Step 2: shift the right-hand mini table (G:K) up by one row.
# Old layout:            New layout:
#   row3 G:I header         (merged into row2 above)
#   row4 G:I data   ----->  row3 G:I data
#   row5 G:I data   ----->  row4 G:I data
#   row6 H:K diff   ----->  row5 H:K diff (formulas now reference row3/row4)
# Stage the formats we need in a scratch area (row 50) first, since source and
# destination ranges overlap (row4/row5 are both read from and written to).
# ---------------------------------------------------------------------------
CopyFormat "G4" "G50"
CopyFormat "H4" "H50"
CopyFormat "I4" "I50"
CopyFormat "G5" "G51"
CopyFormat "H5" "H51"
CopyFormat "I5" "I51"
CopyFormat "H6" "H52"
CopyFormat "I6" "I52"
CopyFormat "J6" "J52"
CopyFormat "K6" "K52"

# now remove the old footprint for rows 4-6 completely (value + style), since
# none of those cells should remain in the output.
$ws.Range("G4:K6").Clear()

# bring the staged formats down into their new homes, then set the shifted
# values/formulas
CopyFormat "G50" "G3"
CopyFormat "H50" "H3"
CopyFormat "I50" "I3"
CopyFormat "G51" "G4"
CopyFormat "H51" "H4"
CopyFormat "I51" "I4"
CopyFormat "H52" "H5"
CopyFormat "I52" "I5"
CopyFormat "J52" "J5"
CopyFormat "K52" "K5"
$ws.Range("G50:K52").Clear()

$ws.Range("G3").Value = $origG4
$ws.Range("H3").Value = $origH4
$ws.Range("I3").Value = $origI4

$ws.Range("G4").Value = $origG5
$ws.Range("H4").Value = $origH5
$ws.Range("I4").Value = $origI5

$ws.Range("H5").Value = "Difference"
$ws.Range("I5").Formula = "=I3-I4"
$ws.Range("J5").Formula = "=100%+I5/I3-1"
$ws.Range("K5").Value = "<-- % lost NVD, some loss was expected"

# ---------------------------------------------------------------------------
# Step 3: second validation block, rows 15-24 (copy structure of rows 2-7,
# with a new file-comparison day and validation formulas against the first
# block).
# ---------------------------------------------------------------------------

# Row 15: headers. A15:D15 bold "Input" variant, G15:I15 plain bold header
# (same as the very first, pre-edit style of the right-hand table header).
CopyFormat "A2" "A15"
$ws.Range("A15").Value = "DATE"
$ws.Range("A15").Font.Bold = $true
CopyFormat "B2" "B15"
$ws.Range("B15").Value = "FILE"
$ws.Range("B15").Font.Bold = $true
CopyFormat "C2" "C15"
$ws.Range("C15").Value = "Count"
$ws.Range("C15").Font.Bold = $true
CopyFormat "D2" "D15"
$ws.Range("D15").Value = "Difference"
$ws.Range("D15").Font.Bold = $true

CopyFormat "E2" "E15"
CopyFormat "F2" "F15"

CopyFormat "D2" "G15"
$ws.Range("G15").Style = "Normal"
$ws.Range("G15").Font.Bold = $true
$ws.Range("G15").Value = "DATE"
CopyFormat "G15" "H15"
$ws.Range("H15").Value = "FILE"
CopyFormat "G15" "I15"
$ws.Range("I15").Value = "NVD billed"

# Row 16: blank left side (A16:F16 same plain style as E/F col), right side
# new data row.
CopyFormat "E2" "A16"
CopyFormat "E2" "B16"
CopyFormat "E2" "C16"
CopyFormat "E2" "D16"
CopyFormat "E2" "E16"
CopyFormat "E2" "F16"

CopyFormat "G3" "G16"
$ws.Range("G16").Value = 44173
CopyFormat "H3" "H16"
$ws.Range("H16").Value = "Billed NVD"
CopyFormat "I3" "I16"
$ws.Range("I16").Value = 54207062

# Row 17: new "Product Detail" data + validation diff vs. row 3.
CopyFormat "A3" "A17"
$ws.Range("A17").Value = 44173
CopyFormat "B3" "B17"
$ws.Range("B17").Value = "Product Detail"
CopyFormat "C3" "C17"
$ws.Range("C17").Value = 47632684
CopyFormat "D5" "D17"
$ws.Range("D17").Formula = "=C17-C3"

CopyFormat "G4" "G17"
$ws.Range("G17").Value = 44173
CopyFormat "H4" "H17"
$ws.Range("H17").Value = "PD + NVD (clean)"
CopyFormat "I4" "I17"
$ws.Range("I17").Value = 53006603

# Row 18: new "Billed NVD" data + validation diff vs. row 4, plus the
# Difference/percent-lost formulas for the new block's right-hand table.
CopyFormat "A4" "A18"
$ws.Range("A18").Value = 44173
CopyFormat "B4" "B18"
$ws.Range("B18").Value = "Billed NVD"
CopyFormat "C4" "C18"
$ws.Range("C18").Value = 4212233
CopyFormat "D6" "D18"
$ws.Range("D18").Formula = "=C18-C4"

CopyFormat "H5" "H18"
$ws.Range("H18").Value = "Difference"
CopyFormat "I5" "I18"
$ws.Range("I18").Formula = "=I16-I17"
CopyFormat "J5" "J18"
$ws.Range("J18").Formula = "=100%+I18/I16-1"
CopyFormat "K5" "K18"
$ws.Range("K18").Value = "<-- % lost NVD, some loss was expected"

# Row 19: "Product Detail (clean) - new" + validation diff.
CopyFormat "A5" "A19"
$ws.Range("A19").Value = 44173
CopyFormat "B5" "B19"
$ws.Range("B19").Value = "Product Detail (clean) - new"
CopyFormat "C5" "C19"
$ws.Range("C19").Value = 47632684
CopyFormat "D5" "D19"
$ws.Range("D19").Formula = "=C17-C19"
CopyFormat "E5" "E19"
CopyFormat "F5" "F19"

# Row 20: "Billed NVD (clean) - new" + validation diff.
CopyFormat "A6" "A20"
$ws.Range("A20").Value = 44173
CopyFormat "B6" "B20"
$ws.Range("B20").Value = "Billed NVD (clean) - new"
CopyFormat "C6" "C20"
$ws.Range("C20").Value = 4212233
CopyFormat "D6" "D20"
$ws.Range("D20").Formula = "=C18-C20"
CopyFormat "E6" "E20"
CopyFormat "F6" "F20"

# Row 21: "Prod Detail and Billed NVD (clean) - new" + validation diff.
CopyFormat "A7" "A21"
$ws.Range("A21").Value = 44173
CopyFormat "B7" "B21"
$ws.Range("B21").Value = "Prod Detail and Billed NVD (clean) - new"
CopyFormat "C7" "C21"
$ws.Range("C21").Value = 14651842
CopyFormat "D7" "D21"
$ws.Range("D21").Formula = "=C21-C7"
CopyFormat "E7" "E21"
CopyFormat "F7" "F21"
CopyFormat "A7" "G21"
$ws.Range("G21").ClearContents()

# Rows 22-24: trailing blank rows with the date/count formatting carried
# down the columns.
CopyFormat "A7" "A22"
$ws.Range("A22").ClearContents()
CopyFormat "C7" "C22"
$ws.Range("C22").ClearContents()

CopyFormat "A7" "A23"
$ws.Range("A23").ClearContents()
CopyFormat "C7" "C23"
$ws.Range("C23").ClearContents()

CopyFormat "A7" "A24"
$ws.Range("A24").ClearContents()
CopyFormat "C7" "C24"
$ws.Range("C24").ClearContents()

# ---------------------------------------------------------------------------
# Step 4: column widths for the new columns D and K.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 10.42578125
$ws.Columns("K").ColumnWidth = 36.7109375

# Remove the old selection on G6 so the saved view just shows the default.
$ws.Range("A1").Select() | Out-Null
